$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old header row ("municipios", "CONFIRMADOS", "ÓBITOS") that lived
# in row 2. This shifts every data row up by one, so the first data row
# ("aguas de lindoia") now becomes row 2, and the last row (previously 109)
# becomes row 108.
$ws.Rows(2).Delete()
